$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 with the same bold header style used by A1:E1,
# plus the two trailing empty-but-styled cells (G1/H1) that Excel leaves
# behind after the header row formatting was extended.
$ws.Range("F1").Value = "geographic_location"
$ws.Range("F1:H1").Font.Bold = $true

# Data rows 2-70 (airport_id 1-69) -> Europe
$ws.Range("F2:F70").Value = "Europe"

# Data rows 71-102 (airport_id 70-101) -> North America
$ws.Range("F71:F102").Value = "North America"

# Leave the selection on F1 (also drops the old topLeftCell scroll anchor)
[void]$ws.Range("F1").Select()
